$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix invalid facility utilisation values so every row uses the same facility id
# Row 5 (facility E5/F5/G5): facility limit + utilisation values were shifted/incorrect
$ws.Range("E5").Value = 600000
$ws.Range("F5").Value = 761579.37
$ws.Range("G5").Value = 3938753.8

# Row 6 (facility E6/F6/G6): facility limit + utilisation currency value corrections
$ws.Range("E6").Value = 600000
$ws.Range("G6").Value = 761579.37

# Column E:G now share a uniform (best-fit) width since the values line up
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 15.5

# Update the active selection to reflect the edited range
$ws.Range("E5:H6").Select()
